# Insert a new weekly price record above row 29 (Feria Lagunitas de Puerto
# Montt - Haba), pushing the existing rows 29..69 down to 30..70. The new
# row reuses the old row-29 record but with an updated date and price.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 29..69 down to 30..70, leaving an empty row 29 behind.
$ws.Rows.Item(29).Insert()

# Populate the newly inserted row 29 with the new data point.
$ws.Range("A29").Value = 4
$ws.Range("B29").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C29").Value = "Los Lagos"
$ws.Range("D29").Value = 44533
$ws.Range("E29").Value = 10
$ws.Range("F29").Value = 100112026
$ws.Range("G29").Value = "Haba"
$ws.Range("H29").Value = "Sin especificar"
$ws.Range("I29").Value = "Primera"
$ws.Range("J29").Value = 120
$ws.Range("K29").Value = 14000
$ws.Range("L29").Value = 14000
$ws.Range("M29").Value = 14000
$ws.Range("N29").Value = "`$/saco 25 kilos"
$ws.Range("O29").Value = "Región del Maule"
$ws.Range("P29").Value = 560
$ws.Range("Q29").Value = 25
$ws.Range("R29").Value = "Hortaliza"
